$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row -> new value for column F ("dSF")
$updates = @{
    4  = 4
    5  = -2
    7  = 2
    14 = -2
    15 = -2
    17 = -1
    20 = 1
    21 = 2
    22 = -1
    25 = -3
    26 = 3
    29 = 0
    30 = 2
    32 = -11
    34 = -2
    38 = 4
    54 = 1
    57 = 1
    61 = 5
    66 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
